# Applies the "Updated project name in docs" edit to CSPSubGuide.docx
#
# Summary of changes:
#  1. Title: wrap "Azure" in a _GoBack bookmark (splits the single title
#     run into three runs: "Preparing " | "Azure" | " CSP Subscription").
#  2. Overview paragraph: "...that you can use with Project Phoenix."
#     becomes "...that you can use with Azure Pack Connector." (the
#     trailing run is split into three runs).
#  3. Pre-requisites: "Install Windows Azure " -> "Install Azure ".
#  4. "...register the subscription in WAP:" becomes
#     "...register the subscription in Azure Pack:" (the run is split
#     into three runs).
#
# Note: the original document already has a stray _GoBack bookmark
# sitting at the very end of the document (after the screenshot).
# Bookmark names must be unique, so re-adding a bookmark named
# "_GoBack" around "Azure" in the title (step 1) automatically moves
# the existing one there - no separate delete step is required (and
# Range.Bookmarks.Exists()/() in this host is not range-scoped, so
# trying to look it up through a different range is unreliable anyway).

$d = $word.ActiveDocument

# --- 1. Title: "Preparing Azure CSP Subscription" ----------------------
# Wrap "Azure" (inside the Title paragraph) with a _GoBack bookmark.
# Adding the bookmark naturally splits the run that contained the title
# text into three runs, matching the target markup, without touching
# the surrounding text. Because a bookmark name is unique document-wide,
# this also relocates (rather than duplicates) the pre-existing _GoBack
# bookmark that used to sit at the end of the document.
$titleRng = $d.Paragraphs(1).Range.Duplicate
$titleRng.Find.Execute("Azure", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $titleRng)

# --- 2. Overview paragraph: Project Phoenix -> Azure Pack Connector ----
# Replace "Project Phoenix." (including the trailing period) with
# "Azure Pack Connector." via a temporary bookmark's Range.Text - this
# keeps the neighbouring runs intact and avoids Word re-flagging every
# run in the paragraph with a stray xml:space="preserve".
$overviewRng = $d.Paragraphs(3).Range.Duplicate
$overviewRng.Find.Execute("Project Phoenix.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("TmpOverview", $overviewRng)
$d.Bookmarks("TmpOverview").Range.Text = "Azure Pack Connector."
$d.Bookmarks("TmpOverview").Delete()

# Now split "Azure Pack Connector." into "Azure Pack Connector" | "."
$overviewSplit = $d.Paragraphs(3).Range.Duplicate
$overviewSplit.Find.Execute("Azure Pack Connector", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("TmpOverviewSplit", $overviewSplit)
$d.Bookmarks("TmpOverviewSplit").Delete()

# --- 3. Pre-requisites: Install Windows Azure -> Install Azure ---------
$d.Paragraphs(5).Range.Find.Execute("Install Windows Azure ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Install Azure ", 2)

# --- 4. "...register the subscription in WAP:" -> "...Azure Pack:" -----
$wapRng = $d.Paragraphs(23).Range.Duplicate
$wapRng.Find.Execute("WAP:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("TmpWap", $wapRng)
$d.Bookmarks("TmpWap").Range.Text = "Azure Pack:"
$d.Bookmarks("TmpWap").Delete()

# Now split "Azure Pack:" into "Azure Pack" | ":"
$wapSplit = $d.Paragraphs(23).Range.Duplicate
$wapSplit.Find.Execute("Azure Pack", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("TmpWapSplit", $wapSplit)
$d.Bookmarks("TmpWapSplit").Delete()
